# Update the three monthly "Ingresos" figures that changed in this commit.
# (row 5 = Diciembre 2023, row 6 = Enero 2024, row 7 = Febrero 2024)
# The SUBTOTAL formula in C11 recalculates automatically to reflect the new total.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ingresos")

$ws.Range("C5").Value = 819.0
$ws.Range("C6").Value = 3245.0
$ws.Range("C7").Value = 369.0
